$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# --- Add the new log entry (row 12) for the Smallest/Second-Smallest problem ---
# (populate F12's URL text first so the new shared strings are interned in the
#  same order as the target workbook: URL = index 72, then the corrected C7
#  caption = index 73)
$ws2.Range("A12").Value = 8
$ws2.Range("B12").Value = 45818
$ws2.Range("C12").Value = "Find the smallest and second smallest element in an array"
$ws2.Range("D12").Value = 1
$ws2.Range("E12").Value = "Easy"
$ws2.Range("F12").Value = "https://www.geeksforgeeks.org/to-find-smallest-and-second-smallest-element-in-an-array/"

# --- Fix the mislabeled entry at C7 (it actually documents "second largest",
#     not "smallest and second smallest") ---
$ws2.Range("C7").Value = "Find the second larges element in an array"

$ws2.Range("A12:E12").HorizontalAlignment = -4131
$ws2.Range("B12").NumberFormat = "m/d/yyyy"

$ws2.Hyperlinks.Add($ws2.Range("F12"), "https://www.geeksforgeeks.org/to-find-smallest-and-second-smallest-element-in-an-array/") | Out-Null

# --- Move the active tab / selection from Sheet1 to Sheet2 ---
$ws2.Activate() | Out-Null
$ws2.Range("F15").Select() | Out-Null
